{"js": "// Separate the \"Target Platform Analysis\" reference from the Project\n// Outline reference in the Document Purpose paragraph, bump the\n// \"2023 Group Project\" citation number from [2] to [3], and tweak the\n// wording of the first objective bullet.\n\n// --- Change 1: Document Purpose paragraph -------------------------------\nlet results = context.document.body.search(\n  \"Project Outline[1] document for correct context.\",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"Project Outline[1] document and Target Platform Analysis document[2] for correct context.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// --- Change 2: \"requirements specification\" paragraph citation number ---\nresults = context.document.body.search(\n  \"2023 Group Project[2]. \",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"2023 Group Project[3]. \", \"Replace\");\n  await context.sync();\n}\n\n// --- Change 3: first bullet under \"The main objectives of this document\" -\nresults = context.document.body.search(\n  \"To define all requirements criteria of the project and what is expected in a final product\",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"To define all requirement criteria of the project and what is expected in the final product\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# Separate the \"Target Platform Analysis\" reference from the Project\n# Outline reference in the Document Purpose paragraph, bump the\n# \"2023 Group Project\" citation number from [2] to [3], and tweak the\n# wording of the first objective bullet.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: Document Purpose paragraph -------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Project Outline[1] document for correct context.\"\nif ($find.Execute()) {\n    $find.Parent.Text = \"Project Outline[1] document and Target Platform Analysis document[2] for correct context.\"\n}\n\n# --- Change 2: \"requirements specification\" paragraph citation number ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"2023 Group Project[2]. \"\nif ($find.Execute()) {\n    $find.Parent.Text = \"2023 Group Project[3]. \"\n}\n\n# --- Change 3: first bullet under \"The main objectives of this document\" -\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"To define all requirements criteria of the project and what is expected in a final product\"\nif ($find.Execute()) {\n    $find.Parent.Text = \"To define all requirement criteria of the project and what is expected in the final product\"\n}\n"}
